$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 empty rows starting at row 9, pushing the existing data
# (old rows 9-12) down to rows 14-17.
$ws.Range("A9:A13").EntireRow.Insert()

# Touch row 13 (a no-op write of its already-default outline level) so the
# freshly inserted, still-empty row gets materialized as its own row
# element in the sheet XML instead of being silently skipped.
$ws.Rows("13").OutlineLevel = 0
